# Regenerate s_val data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and the derived sum column G
# for rows 2-9 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 16.98373111632243;  E = 0.4998867070740569; G = 22.31973251085698 }
    3 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    4 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    5 = @{ B = 3.182878228561681; C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 4.733082622659194 }
    6 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    7 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 16.98373111632243;  E = 0.4998867070740569; G = 22.31973251085698 }
    8 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 3.082599426703578; E = 0.4998867070740569; G = 8.418600821238126 }
    9 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
